$d = $word.ActiveDocument

# The heading paragraph reads "Veranstaltungen des KV München 2023" and must
# become "...2024"; the floating _GoBack bookmark (Word's "last edit location"
# marker) must move from its old spot (after "31.12") to right after the new
# "24".

# 1) Locate the run that holds the last two digits of the year in the heading,
#    scoped to start right after "KV München " so we only touch that run
#    (there are other, unrelated "23" runs earlier in the document).
$probe = $d.Content
$null = $probe.Find.Execute("KV München ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$afterMuenchen = $probe.End

$yearRange = $d.Range($afterMuenchen, $d.Content.End)
$null = $yearRange.Find.Execute("23", $false, $false, $true, $false, $false, $true, 1, $false, "24X", 2)

# 2) Re-plant the _GoBack bookmark immediately after the new "24" (collapsed,
#    zero-length range) while the trailing placeholder "X" still keeps that
#    spot from being the very end of the paragraph; Word only ever keeps a
#    single _GoBack bookmark, so adding the new one removes the stale one
#    that used to sit after "31.12".
$xPos = $yearRange.End - 1
$bmRange = $d.Range($xPos, $xPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# 3) Drop the placeholder now that the bookmark is anchored in place.
$xRange = $d.Range($xPos, $xPos + 1)
$xRange.Text = ""
